$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.581.48'
$ws.Range('E2').Value = '  -1.95%  '
$ws.Range('D3').Value = '3.537.97'
$ws.Range('E3').Value = '  -3.62%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').Value = "'189.21"
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('D6').Value = "'563.21"
$ws.Range('E6').Value = '  -5.99%  '
$ws.Range('D7').Value = '3.532.56'
$ws.Range('E7').Value = '  -3.68%  '
$ws.Range('D8').Value = "'0.608"
$ws.Range('E8').Value = '  -2.32%  '
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').Value = "'0.666"
$ws.Range('E10').Value = '  -5.95%  '
$ws.Range('D11').Value = "'54.94"
$ws.Range('E11').Value = '  -5.93%  '
$ws.Range('D12').Value = "'0.147"
$ws.Range('E12').Value = '  -4.92%  '
$ws.Range('D13').Value = "'0.0000264"
$ws.Range('E13').Value = '  -4.32%  '
$ws.Range('D14').Value = "'9.70"
$ws.Range('E14').Value = '  -5.19%  '
$ws.Range('D15').Value = '4.116.28'
$ws.Range('E15').Value = '  -3.59%  '
$ws.Range('D16').Value = '3.559.30'
$ws.Range('E16').Value = '  -3.40%  '
$ws.Range('E17').Value = '  -1.49%  '
$ws.Range('D18').Value = '66.522.21'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('D19').Value = "'17.97"
$ws.Range('E19').Value = '  -5.51%  '
$ws.Range('D20').Value = "'12.00"
$ws.Range('E20').Value = '  -4.58%  '
$ws.Range('D21').Value = '1.05'
$ws.Range('E21').Value = '  -6.68%  '
$ws.Range('D22').Value = "'395.05"
$ws.Range('E22').Value = '  -1.59%  '
$ws.Range('D23').Value = "'4.10"
$ws.Range('E23').Value = '  -8.01%  '
$ws.Range('D24').Value = "'84.94"
$ws.Range('E24').Value = '  -3.78%  '
$ws.Range('D25').Value = "'11.29"
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').Value = "'2.88"
$ws.Range('E26').Value = '  -3.05%  '
$ws.Range('D27').Value = '12.26'
$ws.Range('E27').Value = '  -2.62%  '
$ws.Range('E28').Value = '  +0.72%  '
$ws.Range('D29').Value = "'3.58"
$ws.Range('E29').Value = '  -3.46%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'8.80"
$ws.Range('E30').Value = '  -5.98%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').Value = "'7.58"
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').Value = "'30.80"
$ws.Range('E32').Value = '  -3.64%  '
$ws.Range('D33').Value = "'634.00"
$ws.Range('E33').Value = '  +2.70%  '
$ws.Range('D34').Value = '11.92'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('D35').Value = "'63.49"
$ws.Range('E35').Value = '  -5.15%  '
$ws.Range('D36').Value = '0.112'
$ws.Range('E36').Value = '  -4.92%  '
$ws.Range('D37').Value = "'41.48"
$ws.Range('E37').Value = '  -9.35%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = "'1.00"
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = "'0.396"
$ws.Range('E39').Value = '  -1.10%  '
$ws.Range('D40').Value = '0.0₃0748'
$ws.Range('E40').Value = '  -4.72%  '
$ws.Range('D41').Value = '3.173.47'
$ws.Range('E41').Value = '  +11.38%  '
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = "'2.90"
$ws.Range('E44').Value = '  -0.83%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = "'2.62"
$ws.Range('E45').Value = '  +1.88%  '
$ws.Range('D46').Value = "'0.0407"
$ws.Range('E46').Value = '  -5.21%  '
$ws.Range('E47').Value = '  -4.23%  '
$ws.Range('E48').Value = '  -5.79%  '
$ws.Range('D49').Value = "'140.73"
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').Value = "'2.51"
$ws.Range('E50').Value = '  -5.38%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = "'8.33"
$ws.Range('E51').Value = '  -7.32%  '
